$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 190
$ws1.Range("F5").Value = 3443
$ws1.Range("F6").Value = 352
$ws1.Range("F8").Value = 426

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 190
$ws4.Range("F5").Value = 3443
$ws4.Range("F6").Value = 352
$ws4.Range("F10").Value = 426
